# Apply the commit: "Looked at altitude of 300 km and 5 deg elevation"
#
# Underlying model: the "Input" worksheet holds named input cells that the
# "UHF" and "S-Band" worksheets derive their link-budget calculations from.
# Changing the inputs below is sufficient for Excel to recalculate every
# dependent formula cell on the UHF and S-Band sheets automatically.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")

# S/C Altitude (alt, Input!C5): 500,000 m -> 300,000 m (300 km)
$wsInput.Range("C5").Value = 300000

# Elevation Angle (del, Input!C6): 20 deg -> 5 deg
$wsInput.Range("C6").Value = 5

# Estimated GS Losses (GSL, Input!C16): 2 dB -> 3.6 dB
$wsInput.Range("C16").Value = 3.6

# Estimated S/C Losses (SCL, Input!C20): 2.2 dB -> 2 dB
$wsInput.Range("C20").Value = 2

# Atmospheric Losses (AIL, Input!C24): 1.1 dB -> 2.1 dB
$wsInput.Range("C24").Value = 2.1

# Force a full recalculation so every formula cell (UHF / S-Band sheets)
# picks up the new inputs before the workbook is saved.
$excel.CalculateFull()

# Update the selections recorded in each sheet view to match where the
# author last clicked while reviewing the change.
$wsInput.Activate()
$wsInput.Range("C25").Select()

$wsUHF = $wb.Worksheets.Item("UHF")
$wsUHF.Activate()
$wsUHF.Range("A10").Select()

$wsSBand = $wb.Worksheets.Item("S-Band")
$wsSBand.Activate()
$wsSBand.Range("G24").Select()
